# Apply updated TPM values to sheet1, then remove the last data row
# (Target cluster "Resolving-Mac" / row 5) which is no longer present
# in the recomputed output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value2 = 0.07205133333333334
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 34.97976933333334
$ws.Range("N2").Value2 = 104.939308
$ws.Range("O2").Value2 = 0.4352965780925344
$ws.Range("P2").Value2 = 0.4352965780925344
$ws.Range("Q2").Value2 = 2.520339020159112
$ws.Range("R2").Value2 = 22.683051181432
$ws.Range("S2").Value2 = 0.4352965780925344
$ws.Range("T2").Value2 = 0.4352965780925344

# --- Row 3 updates ---
$ws.Range("G3").Value2 = 0.07205133333333334
$ws.Range("N3").Value2 = 61.03014900000001
$ws.Range("O3").Value2 = 0.2531579017099818
$ws.Range("P3").Value2 = 0.2531579017099818
$ws.Range("Q3").Value2 = 1.465767869660667
$ws.Range("S3").Value2 = 0.2531579017099818
$ws.Range("T3").Value2 = 0.2531579017099818

# --- Row 4 updates ---
$ws.Range("G4").Value2 = 0.07205133333333334
$ws.Range("M4").Value2 = 25.035323
$ws.Range("N4").Value2 = 75.105969
$ws.Range("O4").Value2 = 0.3115455201974837
$ws.Range("P4").Value2 = 0.3115455201974837
$ws.Range("Q4").Value2 = 1.803828402580667
$ws.Range("R4").Value2 = 16.234455623226
$ws.Range("S4").Value2 = 0.3115455201974837
$ws.Range("T4").Value2 = 0.3115455201974837

# --- Remove row 5 (Resolving-Mac target cluster) entirely ---
$ws.Rows.Item(5).Delete()
